$d = $word.ActiveDocument

# --- Step 1: split the run of 102 'a' characters into "A" + 101 'a's ---
# Replace the first character's text ("a" -> "A"); this initially re-merges
# into a single run, so we then nudge the formatting of just that first
# character (set a direct property on and back off) to force Word to keep
# it as a distinct run from the rest of the text, matching the two-run
# structure produced when a single character is retyped in place.
$firstChar = $d.Range(0, 1)
$firstChar.Text = "A"

$firstCharAgain = $d.Range(0, 1)
$firstCharAgain.Bold = $true
$firstCharAgain.Bold = $false

# --- Step 2: add two blank paragraphs, then a paragraph of 'c' characters ---
$tail = $d.Content
$tail.Collapse(0)
$tail.InsertParagraphAfter()

$tail = $d.Content
$tail.Collapse(0)
$tail.InsertParagraphAfter()

$tail = $d.Content
$tail.Collapse(0)
$tail.InsertParagraphAfter()

$tail = $d.Content
$tail.Collapse(0)
$tail.InsertAfter("ccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccccc")
